$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.692.48"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.584.10"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'207.27"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'22.35"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").Value = "'0.254"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("E10").Value = "  -3.36%  "
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "1.808.19"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "1.583.68"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "27.648.42"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "'62.96"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "'217.67"
$ws.Range("E18").Value = "  -4.85%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0694"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  -4.08%  "
$ws.Range("D23").Value = "'9.56"
$ws.Range("E23").Value = "  -4.60%  "
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("D25").Value = "'153.66"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'6.70"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("E29").Value = "  -4.33%  "
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "'0.0464"
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("D33").Value = "1.376.72"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("E36").Value = "  -4.90%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("D40").Value = "'0.817"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.976"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").Value = "'1.79"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D45").Value = "'63.67"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").Value = "'5.23"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "1.719.45"
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("D48").Value = "'87.70"
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("E50").Value = "  -4.91%  "
$ws.Range("E51").Value = "  -1.59%  "
